# Applies the cryptos-list update described in the commit:
# "Updated cryptos list on Thu Mar  7 17:42:40 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'67.634.88"
$ws.Range("E2").Value = "  +1.24%  "

# Row 3
$ws.Range("D3").Value = "'3.860.72"
$ws.Range("E3").Value = "  +1.42%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "'456.44"
$ws.Range("E5").Value = "  +8.28%  "

# Row 6
$ws.Range("D6").Value = "'146.86"
$ws.Range("E6").Value = "  +13.15%  "

# Row 7
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  +2.74%  "

# Row 8
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.746"
$ws.Range("E9").Value = "  +3.68%  "

# Row 10
$ws.Range("E10").Value = "  -3.00%  "

# Row 11
$ws.Range("D11").Value = "'0.0000321"
$ws.Range("E11").Value = "  -7.40%  "

# Row 12
$ws.Range("D12").Value = "'43.88"
$ws.Range("E12").Value = "  +8.12%  "

# Row 13
$ws.Range("D13").Value = "'10.37"
$ws.Range("E13").Value = "  +2.47%  "

# Row 14
$ws.Range("D14").Value = "'4.464.19"
$ws.Range("E14").Value = "  +1.40%  "

# Row 15
$ws.Range("D15").Value = "'14.83"
$ws.Range("E15").Value = "  -4.43%  "

# Row 16
$ws.Range("D16").Value = "'3.914.52"
$ws.Range("E16").Value = "  +2.65%  "

# Row 17
$ws.Range("E17").Value = "  -0.22%  "

# Row 18
$ws.Range("D18").Value = "'20.14"
$ws.Range("E18").Value = "  +2.55%  "

# Row 19
$ws.Range("E19").Value = "  +8.41%  "

# Row 20
$ws.Range("D20").Value = "'67.739.66"
$ws.Range("E20").Value = "  +1.07%  "

# Row 21
$ws.Range("D21").Value = "'429.00"
$ws.Range("E21").Value = "  +5.71%  "

# Row 22
$ws.Range("D22").Value = "'14.84"
$ws.Range("E22").Value = "  -2.85%  "

# Row 23
$ws.Range("D23").Value = "'3.25"
$ws.Range("E23").Value = "  +7.11%  "

# Row 24
$ws.Range("D24").Value = "'86.80"
$ws.Range("E24").Value = "  +3.63%  "

# Row 25
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "'10.44"
$ws.Range("E25").Value = "  +15.78%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'3.51"
$ws.Range("E26").Value = "  +9.37%  "

# Row 27
$ws.Range("D27").Value = "'37.53"
$ws.Range("E27").Value = "  +1.27%  "

# Row 28
$ws.Range("D28").Value = "'9.85"
$ws.Range("E28").Value = "  -1.15%  "

# Row 29
$ws.Range("D29").Value = "'5.50"
$ws.Range("E29").Value = "  +1.06%  "

# Row 30
$ws.Range("D30").Value = "'740.67"
$ws.Range("E30").Value = "  +2.89%  "

# Row 31
$ws.Range("D31").Value = "'13.78"
$ws.Range("E31").Value = "  +8.92%  "

# Row 32
$ws.Range("E32").Value = "  +11.49%  "

# Row 33
$ws.Range("E33").Value = "  -1.09%  "

# Row 34
$ws.Range("D34").Value = "'43.39"
$ws.Range("E34").Value = "  +12.64%  "

# Row 35
$ws.Range("E35").Value = "  +6.49%  "

# Row 36
$ws.Range("D36").Value = "'57.34"
$ws.Range("E36").Value = "  +3.99%  "

# Row 37
$ws.Range("E37").Value = "  +3.15%  "

# Row 38
$ws.Range("E38").Value = "  +0.07%  "

# Row 39
$ws.Range("D39").Value = "'0.0477"
$ws.Range("E39").Value = "  +5.45%  "

# Row 40
$ws.Range("D40").Value = "'0.359"
$ws.Range("E40").Value = "  +15.98%  "

# Row 41
$ws.Range("D41").Value = "'2.96"
$ws.Range("E41").Value = "  +1.76%  "

# Row 42
$ws.Range("D42").Value = "'2.68"
$ws.Range("E42").Value = "  +20.59%  "

# Row 43
$ws.Range("E43").Value = "  -7.89%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.22%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.140"
$ws.Range("E45").Value = "  +4.61%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.32"
$ws.Range("E46").Value = "  +6.91%  "

# Row 47
$ws.Range("E47").Value = "  +3.99%  "

# Row 48
$ws.Range("D48").Value = "'2.15"
$ws.Range("E48").Value = "  +5.63%  "

# Row 49
$ws.Range("E49").Value = "  +6.21%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'144.49"
$ws.Range("E50").Value = "  +0.85%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'2.88"
$ws.Range("E51").Value = "  +3.22%  "
